$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 392; this shifts the existing rows
# 392-494 down to 394-496 (preserving all of their data), matching the
# diff where every pair of "Primera"/"Segunda" rows from 392 onward is
# pushed down by 2 rows and the final two pairs reappear as new rows
# 495-496.
$ws.Rows("392:393").Insert()

# Populate the two brand-new rows (392 and 393) with the new weekly
# price-report entries for Betarraga / Femacal de La Calera.
$ws.Range("A392").Value = 3
$ws.Range("B392").Value = "Femacal de La Calera"
$ws.Range("C392").Value = "Coquimbo"
$ws.Range("D392").Value = 44551
$ws.Range("E392").Value = 5
$ws.Range("F392").Value = 100114014
$ws.Range("G392").Value = "Betarraga"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 180
$ws.Range("K392").Value = 500
$ws.Range("L392").Value = 500
$ws.Range("M392").Value = 500
$ws.Range("N392").Value = "$/paquete 4 unidades"
$ws.Range("O392").Value = "Provincia de Quillota"
$ws.Range("P392").Value = 125
$ws.Range("Q392").Value = 4
$ws.Range("R392").Value = "Hortaliza"

$ws.Range("A393").Value = 3
$ws.Range("B393").Value = "Femacal de La Calera"
$ws.Range("C393").Value = "Coquimbo"
$ws.Range("D393").Value = 44551
$ws.Range("E393").Value = 5
$ws.Range("F393").Value = 100114014
$ws.Range("G393").Value = "Betarraga"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Segunda"
$ws.Range("J393").Value = 160
$ws.Range("K393").Value = 400
$ws.Range("L393").Value = 400
$ws.Range("M393").Value = 400
$ws.Range("N393").Value = "$/paquete 4 unidades"
$ws.Range("O393").Value = "Provincia de Quillota"
$ws.Range("P393").Value = 100
$ws.Range("Q393").Value = 4
$ws.Range("R393").Value = "Hortaliza"

# Make sure the date cells keep the date number format used by the
# rest of column D.
$ws.Range("D392:D393").NumberFormat = $ws.Range("D394").NumberFormat

$ws.Range("A1").Select()
